{"js": "// Replace each \"divide\" expression in the worksheet table with the new one.\n// Every \"from\" string occurs exactly once in the document, so a plain\n// text search + full replace is safe and unambiguous.\nconst replacements = [\n  [\"600\u00f75=120, 0\", \"139\u00f76=23, 1\"],\n  [\"755\u00f78=94, 3\", \"538\u00f74=134, 2\"],\n  [\"692\u00f74=173, 0\", \"200\u00f73=66, 2\"],\n  [\"883\u00f78=110, 3\", \"117\u00f78=14, 5\"],\n  [\"957\u00f73=319, 0\", \"477\u00f73=159, 0\"],\n  [\"157\u00f77=22, 3\", \"443\u00f77=63, 2\"],\n  [\"829\u00f73=276, 1\", \"282\u00f79=31, 3\"],\n  [\"685\u00f79=76, 1\", \"958\u00f78=119, 6\"],\n  [\"601\u00f73=200, 1\", \"564\u00f76=94, 0\"],\n  [\"737\u00f79=81, 8\", \"860\u00f73=286, 2\"],\n  [\"940\u00f73=313, 1\", \"288\u00f76=48, 0\"],\n  [\"323\u00f72=161, 1\", \"193\u00f76=32, 1\"],\n  [\"905\u00f72=452, 1\", \"258\u00f77=36, 6\"],\n  [\"599\u00f76=99, 5\", \"814\u00f79=90, 4\"],\n  [\"580\u00f73=193, 1\", \"470\u00f73=156, 2\"],\n  [\"753\u00f75=150, 3\", \"388\u00f72=194, 0\"],\n  [\"474\u00f78=59, 2\", \"113\u00f74=28, 1\"],\n  [\"533\u00f76=88, 5\", \"397\u00f76=66, 1\"],\n  [\"324\u00f79=36, 0\", \"675\u00f74=168, 3\"],\n  [\"172\u00f72=86, 0\", \"905\u00f74=226, 1\"],\n  [\"776\u00f77=110, 6\", \"937\u00f76=156, 1\"],\n  [\"433\u00f76=72, 1\", \"216\u00f74=54, 0\"],\n  [\"606\u00f79=67, 3\", \"978\u00f75=195, 3\"],\n  [\"197\u00f79=21, 8\", \"805\u00f78=100, 5\"],\n  [\"453\u00f79=50, 3\", \"910\u00f78=113, 6\"],\n];\n\nconst body = context.document.body;\n\nfor (const [from, to] of replacements) {\n  const found = body.search(from, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${from}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"divide\" expression in the worksheet table with the new one.\n# Every \"from\" string occurs exactly once in the document, so Find/Replace\n# across the whole document content is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"600\u00f75=120, 0\", \"139\u00f76=23, 1\"),\n    @(\"755\u00f78=94, 3\", \"538\u00f74=134, 2\"),\n    @(\"692\u00f74=173, 0\", \"200\u00f73=66, 2\"),\n    @(\"883\u00f78=110, 3\", \"117\u00f78=14, 5\"),\n    @(\"957\u00f73=319, 0\", \"477\u00f73=159, 0\"),\n    @(\"157\u00f77=22, 3\", \"443\u00f77=63, 2\"),\n    @(\"829\u00f73=276, 1\", \"282\u00f79=31, 3\"),\n    @(\"685\u00f79=76, 1\", \"958\u00f78=119, 6\"),\n    @(\"601\u00f73=200, 1\", \"564\u00f76=94, 0\"),\n    @(\"737\u00f79=81, 8\", \"860\u00f73=286, 2\"),\n    @(\"940\u00f73=313, 1\", \"288\u00f76=48, 0\"),\n    @(\"323\u00f72=161, 1\", \"193\u00f76=32, 1\"),\n    @(\"905\u00f72=452, 1\", \"258\u00f77=36, 6\"),\n    @(\"599\u00f76=99, 5\", \"814\u00f79=90, 4\"),\n    @(\"580\u00f73=193, 1\", \"470\u00f73=156, 2\"),\n    @(\"753\u00f75=150, 3\", \"388\u00f72=194, 0\"),\n    @(\"474\u00f78=59, 2\", \"113\u00f74=28, 1\"),\n    @(\"533\u00f76=88, 5\", \"397\u00f76=66, 1\"),\n    @(\"324\u00f79=36, 0\", \"675\u00f74=168, 3\"),\n    @(\"172\u00f72=86, 0\", \"905\u00f74=226, 1\"),\n    @(\"776\u00f77=110, 6\", \"937\u00f76=156, 1\"),\n    @(\"433\u00f76=72, 1\", \"216\u00f74=54, 0\"),\n    @(\"606\u00f79=67, 3\", \"978\u00f75=195, 3\"),\n    @(\"197\u00f79=21, 8\", \"805\u00f78=100, 5\"),\n    @(\"453\u00f79=50, 3\", \"910\u00f78=113, 6\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
